$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 0.5344385436970592
$ws.Range("E10").Value = 0.6201980835478056
$ws.Range("G10").Value = 6.341402055340052
$ws.Range("H10").Value = 39.51851520369351
$ws.Range("I10").Value = 54.14008274096643
$ws.Range("J10").Value = 45.85991725903357
$ws.Range("K10").Value = 54.14008274096643
$ws.Range("B11").Value = 1.124014212475507
$ws.Range("E11").Value = 3.028873523761085
$ws.Range("G11").Value = 25.99140387354132
$ws.Range("H11").Value = 44.04745699154409
$ws.Range("I11").Value = 29.9611391349146
$ws.Range("J11").Value = 70.03886086508541
$ws.Range("K11").Value = 29.9611391349146
$ws.Range("C12").Value = 3.699075689035789
$ws.Range("E12").Value = 4.517698030842167
$ws.Range("G12").Value = 14.93902192521486
$ws.Range("H12").Value = 67.50435457160448
$ws.Range("I12").Value = 17.55662350318067
$ws.Range("J12").Value = 82.44337649681934
$ws.Range("K12").Value = 17.55662350318067
$ws.Range("D15").Value = 4.588987080161871
$ws.Range("F15").Value = 4.588987080161871
$ws.Range("G15").Value = 22.93056935485755
$ws.Range("H15").Value = 26.57459035111507
$ws.Range("I15").Value = 50.49484029402739
$ws.Range("J15").Value = 49.50515970597262
$ws.Range("K15").Value = 50.49484029402739
$ws.Range("B16").Value = 0.5731952950638293
$ws.Range("E16").Value = 1.080435286159236
$ws.Range("G16").Value = 43.27108330943184
$ws.Range("H16").Value = 38.29205176940727
$ws.Range("I16").Value = 18.43686492116091
$ws.Range("J16").Value = 81.56313507883911
$ws.Range("K16").Value = 18.43686492116091
$ws.Range("C22").Value = 0.02118677
$ws.Range("D22").Value = 0.00776831
$ws.Range("E22").Value = 0.04865666
$ws.Range("F22").Value = 0.00776831
$ws.Range("G22").Value = 48.68392486517937
$ws.Range("H22").Value = 37.54857113791996
$ws.Range("I22").Value = 13.76750399690066
$ws.Range("J22").Value = 86.23249600309933
$ws.Range("K22").Value = 13.76750399690066
$ws.Range("B24").Value = 18.50934330204852
$ws.Range("C24").Value = 33.69065111221759
$ws.Range("D24").Value = 17.59157887518554
$ws.Range("E24").Value = 52.19999441426612
$ws.Range("F24").Value = 17.59157887518554
$ws.Range("G24").Value = 26.52088558783935
$ws.Range("H24").Value = 48.27323632967825
$ws.Range("I24").Value = 25.2058780824824
$ws.Range("J24").Value = 74.7941219175176
$ws.Range("K24").Value = 25.2058780824824
